# BOT; UPDATE DATA
# Adds one new day of COVID-19 stats (2020-04-26, serial 43947) to each of
# the "all", "kobe" and "other" sheets, pushing the trailing footnote
# row(s) down by one row on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all": insert a fresh data row at row 19 (pushes the two footnote
# rows down to 20/21), then fill in the new day's numbers.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")

$wsAll.Rows(19).Insert(-4121)
$wsAll.Range("A19").Value = 43947
$wsAll.Range("B19").Value = 251
$wsAll.Range("C19").Value = 224
$wsAll.Range("D19").Value = 133
$wsAll.Range("E19").Value = 123
$wsAll.Range("F19").Value = 10
$wsAll.Range("G19").Value = 3
$wsAll.Range("H19").Value = 88

$wsAll.Range("C20").Select()

# ---------------------------------------------------------------------
# Sheet "kobe": correct the last existing data row (73), then insert a
# new data row at 74 (pushes the footnote row down to 75) and fill it in.
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

$wsKobe.Range("D73").Value = 6
$wsKobe.Range("E73").Value = 251

$wsKobe.Rows(74).Insert(-4121)
$wsKobe.Range("A74").Value = 43947
$wsKobe.Range("B74").Value = 49
$wsKobe.Range("C74").Value = 1681
$wsKobe.Range("D74").Value = 0
$wsKobe.Range("E74").Value = 251
$wsKobe.Range("F74").Value = 128
$wsKobe.Range("G74").Value = 119
$wsKobe.Range("H74").Value = 9
$wsKobe.Range("I74").Value = 3
$wsKobe.Range("J74").Value = 82

# the footnote row that got pushed to 75 picks up a styled-but-empty A75
# cell (matches the original author's edit) - copy the date format down
# without carrying any value.
$wsKobe.Range("A74").Copy()
$wsKobe.Range("A75").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsKobe.Range("A75").ClearContents()

$wsKobe.Range("A74").Select()

# ---------------------------------------------------------------------
# Sheet "other": row 49 was the footnote placeholder row - turn it into
# a real data row (copying the number formatting down from row 48),
# move the footnote text into row 50 (previously blank), and append a
# fresh blank spacer row 51.
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")

$wsOther.Range("A48:I48").Copy()
$wsOther.Range("A49:I49").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsOther.Range("A49").Value = 43947
$wsOther.Range("B49").Value = 0
$wsOther.Range("C49").Value = 11
$wsOther.Range("D49").Value = 5
$wsOther.Range("E49").Value = 4
$wsOther.Range("F49").Value = 1
$wsOther.Range("G49").Value = 0
$wsOther.Range("H49").Value = 6

$wsOther.Range("B50").Value = "※他自治体において、3月10日以前の感染者の発生はございません。"

$wsOther.Range("A50:I50").Copy()
$wsOther.Range("A51:I51").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsOther.Range("A51:I51").ClearContents()

$wsOther.Range("D50").Select()
